# Add the August 18th, 2020 SSA data row (row 80) to the historical log.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 80

# Column A holds dates that are stored as plain text (shared strings),
# not as real Excel date serials. Temporarily mark the cell as Text so
# typing "2020-08-18" is not auto-converted into a date value, then
# copy the (unstyled) format of an existing data cell in column A back
# onto it so the cell itself ends up with no explicit style, matching
# the rest of the column.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2020-08-18"
$ws.Cells.Item($row, 1).Style = $ws.Cells.Item(2, 1).Style

$ws.Cells.Item($row, 2).Value = 531239
$ws.Cells.Item($row, 3).Value = 584293
$ws.Cells.Item($row, 4).Value = 81175
$ws.Cells.Item($row, 5).Value = 57774
$ws.Cells.Item($row, 6).Value = 26.15
